# Insert a new weekly price-observation row for "Vega Monumental Concepción -
# Plátano" right above the current row 744. This pushes every existing row
# at/after 744 down by one (744->745, ..., 794->795) and extends the used
# range from A1:T794 to A1:T795.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(744).Insert()

# Populate the freshly inserted row with the new observation. The
# market/region/product/category/unit/origin columns are constant across
# every row of this sheet, so they are simply repeated here.
$ws.Cells.Item(744, 1).Value = 11
$ws.Cells.Item(744, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(744, 3).Value = "Bíobío"
$ws.Cells.Item(744, 4).Value = 45106
$ws.Cells.Item(744, 5).Value = 8
$ws.Cells.Item(744, 6).Value = "Fruta"
$ws.Cells.Item(744, 7).Value = 100108
$ws.Cells.Item(744, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(744, 9).Value = 100108006
$ws.Cells.Item(744, 10).Value = "Plátano"
$ws.Cells.Item(744, 11).Value = "Sin especificar"
$ws.Cells.Item(744, 12).Value = "Pintón"
$ws.Cells.Item(744, 13).Value = 1100
$ws.Cells.Item(744, 14).Value = 14000
$ws.Cells.Item(744, 15).Value = 15000
$ws.Cells.Item(744, 16).Value = 14545
$ws.Cells.Item(744, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(744, 18).Value = "Ecuador"
$ws.Cells.Item(744, 19).Value = 727
$ws.Cells.Item(744, 20).Value = 20
